# Updated cryptos list on Wed Jan  3 21:43:05 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.912.86"
$ws.Range("E2").Value = "  -4.71%  "

$ws.Range("D3").Value = "2.225.39"
$ws.Range("E3").Value = "  -5.57%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'317.70"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").Value = "'99.99"
$ws.Range("E6").Value = "  -6.84%  "

$ws.Range("E7").Value = "  -5.97%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.563"
$ws.Range("E9").Value = "  -6.93%  "

$ws.Range("D10").Value = "'37.42"
$ws.Range("E10").Value = "  -7.95%  "

$ws.Range("D11").Value = "'54.11"
$ws.Range("E11").Value = "  -2.67%  "

$ws.Range("D12").Value = "'0.0829"
$ws.Range("E12").Value = "  -9.38%  "

$ws.Range("D13").Value = "'7.80"
$ws.Range("E13").Value = "  -7.20%  "

$ws.Range("E14").Value = "  -2.58%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.563.63"
$ws.Range("E15").Value = "  -5.57%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.863"
$ws.Range("E16").Value = "  -11.34%  "

$ws.Range("D17").Value = "'14.30"
$ws.Range("E17").Value = "  -5.54%  "

$ws.Range("D18").Value = "2.220.67"
$ws.Range("E18").Value = "  -6.04%  "

$ws.Range("D19").Value = "42.816.77"
$ws.Range("E19").Value = "  -4.82%  "

$ws.Range("D20").Value = "'14.94"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("E21").Value = "  -8.54%  "

$ws.Range("D22").Value = "'6.45"
$ws.Range("E22").Value = "  -10.13%  "

$ws.Range("D23").Value = "'65.65"
$ws.Range("E23").Value = "  -9.87%  "

$ws.Range("D24").Value = "'3.16"
$ws.Range("E24").Value = "  -9.69%  "

$ws.Range("D25").Value = "'236.53"
$ws.Range("E25").Value = "  -8.31%  "

$ws.Range("D26").Value = "'2.15"
$ws.Range("E26").Value = "  -6.43%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  -8.89%  "

$ws.Range("E29").Value = "  -4.68%  "

$ws.Range("D30").Value = "'6.39"
$ws.Range("E30").Value = "  -10.61%  "

$ws.Range("D31").Value = "'0.0909"
$ws.Range("E31").Value = "  -5.89%  "

$ws.Range("D32").Value = "'20.50"
$ws.Range("E32").Value = "  -7.77%  "

$ws.Range("D33").Value = "'34.15"
$ws.Range("E33").Value = "  -7.79%  "

$ws.Range("D34").Value = "'156.59"
$ws.Range("E34").Value = "  -6.56%  "

$ws.Range("D35").Value = "'2.77"
$ws.Range("E35").Value = "  -7.01%  "

$ws.Range("D36").Value = "'3.21"
$ws.Range("E36").Value = "  +10.82%  "

$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = "  +12.72%  "

$ws.Range("E38").Value = "  -5.66%  "

$ws.Range("D39").Value = "'3.94"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "'4.47"
$ws.Range("E40").Value = "  -3.97%  "

$ws.Range("E41").Value = "  -8.58%  "

$ws.Range("E42").Value = "  -6.91%  "

$ws.Range("D43").Value = "1.942.66"
$ws.Range("E43").Value = "  +3.36%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "'12.53"
$ws.Range("E45").Value = "  -1.97%  "

$ws.Range("D46").Value = "'89.14"
$ws.Range("E46").Value = "  -10.71%  "

$ws.Range("E47").Value = "  -8.48%  "

$ws.Range("D48").Value = "'5.38"
$ws.Range("E48").Value = "  -4.10%  "

$ws.Range("D49").Value = "'76.52"
$ws.Range("E49").Value = "  -5.75%  "

$ws.Range("D50").Value = "'60.60"

$ws.Range("D51").Value = "'0.868"
$ws.Range("E51").Value = "  +18.95%  "
